# Apply the diff: update dSF (column F) values for several rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = -3
$ws.Range("F6").Value  = 0
$ws.Range("F7").Value  = -3
$ws.Range("F9").Value  = -4
$ws.Range("F11").Value = -2
$ws.Range("F18").Value = 3
$ws.Range("F23").Value = 8
$ws.Range("F27").Value = -4
